# Automatic update of files.
# Rewrites rows 4-11 of the active sheet so that each row's data is
# replaced with the data that (per the diff) now belongs there. This is
# effectively a cyclic re-shuffle of the 8 observation records together
# with a handful of per-record value corrections (coordinates, dates, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 --------------------------------------------------------------
$ws.Range("A4").Value = 74918827
$ws.Range("B4").Value = 103265
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 221144
$ws.Range("F4").Value = "Grönpyrola"
$ws.Range("G4").Value = "Pyrola chlorantha"
$ws.Range("H4").Value = "Sw."
$ws.Range("Q4").Value = 548731.5793738363
$ws.Range("R4").Value = 6316818.087442902

# ---- Row 5 --------------------------------------------------------------
$ws.Range("A5").Value = 74918841
$ws.Range("B5").Value = 73680
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 306
$ws.Range("F5").Value = "Kornig nållav"
$ws.Range("G5").Value = "Chaenotheca chlorella"
$ws.Range("H5").Value = "(Ach.) Müll.Arg."
$ws.Range("Q5").Value = 548718.5936764533
$ws.Range("R5").Value = 6316901.286082684

# ---- Row 6 --------------------------------------------------------------
$ws.Range("A6").Value = 74918848
$ws.Range("B6").Value = 103178
$ws.Range("E6").Value = 221141
$ws.Range("F6").Value = "Gullviva"
$ws.Range("G6").Value = "Primula veris"
$ws.Range("H6").Value = "L."
$ws.Range("Q6").Value = 548759.7599674745
$ws.Range("R6").Value = 6316928.463234497
# Assign the date-like strings while temporarily forcing a text format so
# that Excel does not silently convert them to date serial values.
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2018-11-02"
$ws.Range("Y6").NumberFormat = "General"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2018-11-02"
$ws.Range("AA6").NumberFormat = "General"

# ---- Row 7 --------------------------------------------------------------
$ws.Range("A7").Value = 74918851
$ws.Range("B7").Value = 103265
$ws.Range("E7").Value = 221144
$ws.Range("F7").Value = "Grönpyrola"
$ws.Range("G7").Value = "Pyrola chlorantha"
$ws.Range("H7").Value = "Sw."
$ws.Range("Q7").Value = 548759.7599674745
$ws.Range("R7").Value = 6316928.463234497

# ---- Row 8 --------------------------------------------------------------
$ws.Range("A8").Value = 74918845
$ws.Range("B8").Value = 4711
$ws.Range("E8").Value = 100299
$ws.Range("F8").Value = "Thomsons trägnagare"
$ws.Range("G8").Value = "Cacotemnus thomsoni"
$ws.Range("H8").Value = "(Kraatz, 1881)"
$ws.Range("Q8").Value = 548718.5936764533
$ws.Range("R8").Value = 6316901.286082684

# ---- Row 9 --------------------------------------------------------------
$ws.Range("A9").Value = 74918868
$ws.Range("B9").Value = 73507
$ws.Range("E9").Value = 6428
$ws.Range("F9").Value = "Rostfläck"
$ws.Range("G9").Value = "Arthonia vinosa"
$ws.Range("H9").Value = "Leight."
$ws.Range("L9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("Q9").Value = 548706.4407073758
$ws.Range("R9").Value = 6316820.51631712
$ws.Range("AC9").Value = "På björk och al"

# ---- Row 10 -------------------------------------------------------------
$ws.Range("A10").Value = 74918824
$ws.Range("B10").Value = 88933
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 256335
$ws.Range("F10").Value = "Taggfingersvamp"
$ws.Range("G10").Value = "Ramaria karstenii"
$ws.Range("H10").Value = "(Sacc. & P.Syd.) Corner"
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value = 548761.3979890828
$ws.Range("R10").Value = 6316789.019183607
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "2018-11-09"
$ws.Range("Y10").NumberFormat = "General"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = "2018-11-09"
$ws.Range("AA10").NumberFormat = "General"

# ---- Row 11 -------------------------------------------------------------
$ws.Range("A11").Value = 74918837
$ws.Range("B11").Value = 90638
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 1968
$ws.Range("F11").Value = "Grantaggsvamp"
$ws.Range("G11").Value = "Bankera violascens"
$ws.Range("H11").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("Q11").Value = 548676.2777817947
$ws.Range("R11").Value = 6316878.998678687
$ws.Range("AC11").ClearContents()

# ---- Columns that need to appear as present-but-empty placeholders ------
# Plain assignment of "" does not materialize a cell, so instead copy an
# already-blank cell from the same row onto the target so the column is
# present (but empty), matching the shifted row's original layout.
$ws.Range("I4").Copy($ws.Range("L4"))
$ws.Range("I4").Copy($ws.Range("N4"))
$ws.Range("I7").Copy($ws.Range("L7"))
$ws.Range("I7").Copy($ws.Range("N7"))
$ws.Range("I8").Copy($ws.Range("M8"))
